# Update the NATMI ligand-receptor TPM-derived statistics for Lgi1-Rtn4r
# with newly recomputed TPM-based values (per commit: "update scripts wuth new tpm").
#
# Only the numeric result columns (E..T) on rows 2-5 of the single data sheet
# are affected; identifying columns (A..D) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.6630574638774661
$ws.Range("J2").Value = 0.663057463877466
$ws.Range("O2").Value = 0.858667536176972
$ws.Range("P2").Value = 0.858667536176972
$ws.Range("Q2").Value = 0.02561400552933333
$ws.Range("S2").Value = 0.5693459188514154
$ws.Range("T2").Value = 0.5693459188514154

# Row 3
$ws.Range("I3").Value = 0.6630574638774661
$ws.Range("J3").Value = 0.663057463877466
$ws.Range("S3").Value = 0.09371154502605064
$ws.Range("T3").Value = 0.09371154502605063

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.08283833333333333
$ws.Range("H4").Value = 0.248515
$ws.Range("I4").Value = 0.3369425361225339
$ws.Range("J4").Value = 0.3369425361225339
$ws.Range("O4").Value = 0.858667536176972
$ws.Range("P4").Value = 0.858667536176972
$ws.Range("Q4").Value = 0.01301613880166667
$ws.Range("R4").Value = 0.117145249215
$ws.Range("S4").Value = 0.2893216173255566
$ws.Range("T4").Value = 0.2893216173255566

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.08283833333333333
$ws.Range("H5").Value = 0.248515
$ws.Range("I5").Value = 0.3369425361225339
$ws.Range("J5").Value = 0.3369425361225339
$ws.Range("Q5").Value = 0.002142392589444445
$ws.Range("R5").Value = 0.019281533305
$ws.Range("S5").Value = 0.04762091879697732
$ws.Range("T5").Value = 0.04762091879697731
